# auto increment public id and fix unit tests
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear the "Public ID" column (A) for all data rows - IDs will be
# auto-incremented rather than duplicating the Name text, so the sample
# workbook no longer ships hard-coded values there.
$ws.Range("A2:A11").ClearContents()

# Remove the sample data for Theme 5 - Theme 10 (rows 6-11, columns B:C),
# and restyle those rows (A:C) to match the already-blank placeholder rows
# below them (rows 12-16), since clearing content alone leaves the old
# "text" number format behind.
$ws.Range("B6:C11").ClearContents()
$ws.Range("A6:C11").NumberFormat = "General"

# Drop the now-unused trailing placeholder rows (17-22), shrinking the
# sheet from 22 to 16 rows.
$ws.Range("A17:H22").EntireRow.Delete()
